# Split the "{m:endfor}" run sequence so that the trailing "}" becomes its
# own run positioned AFTER the _GoBack bookmark (instead of being part of
# the "for}" run that currently sits BEFORE the bookmark).
#
# Before:
#   <w:r>{m:end</w:r><w:r>for}</w:r><bookmarkStart/><bookmarkEnd/>
# After:
#   <w:r>{m:end</w:r><w:r>for</w:r><bookmarkStart/><bookmarkEnd/><w:r>}</w:r>

$d = $word.ActiveDocument

# Locate the "{m:endfor}" text (it is split across two runs: "{m:end" and
# "for}") using Find so we do not depend on hard-coded offsets.
$findRange = $d.Content
$found = $findRange.Find.Execute("{m:endfor}", $true, $false, $false, $false, `
                                  $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate '{m:endfor}' in the document"
}

$matchStart = $findRange.Start
$matchEnd = $findRange.End

# The "for}" run is the trailing 4 characters of the match ("for" + "}").
$run2Start = $matchEnd - 4
$run2End = $matchEnd

# 1) Shrink the "for}" run down to "for" in a single atomic replace (using
#    FormattedText so the formatting/rPr is carried over) so the engine
#    does not coalesce it back into the preceding "{m:end" run and does not
#    disturb the bookmark that immediately follows it.
$shrinkSource = $d.Range($run2Start, $run2End - 1)
$shrinkTarget = $d.Range($run2Start, $run2End)
$shrinkTarget.FormattedText = $shrinkSource.FormattedText

# 2) Insert a brand-new run containing "}" right after the bookmark (i.e.
#    right before the paragraph mark) using the same formatting as the
#    surrounding text.
$para = $matchStart
$ownerParagraph = $d.Range($matchStart, $matchStart).Paragraphs(1)
$insertPos = $ownerParagraph.Range.End - 1

$formatSource = $d.Range($matchStart, $matchStart + 1)
$insertRange = $d.Range($insertPos, $insertPos)
$insertRange.FormattedText = $formatSource.FormattedText

$newRunRange = $d.Range($insertPos, $insertPos + 1)
$newRunRange.Text = "}"
